# "massive MDY site template update"
# Insert three new columns (Month, Day, Year) between the existing
# "Transect" column (D) and "Date Sampled" column (old E, now H), splitting
# the sampling date into separate Month/Day/Year fields while keeping the
# original "Date Sampled" column and everything after it intact (shifted
# right by three columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank columns at E:G - this pushes the old E:N (Date Sampled
# .. Notes) block to H:Q, carrying over each cell's style/format, exactly
# like using Excel's "Insert Sheet Columns" on a 3-column selection.
$ws.Range("E1:G1").EntireColumn.Insert()

# New header row labels for the inserted columns.
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# Row 2 -> Transect 1, Date Sampled 8/25/2016 (serial 42607)
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 25
$ws.Range("G2").Value = 2016

# Row 3 -> Transect 2, Date Sampled 8/29/2016 (serial 42611)
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 29
$ws.Range("G3").Value = 2016

# Row 4 -> Transect 3, Date Sampled 8/30/2016 (serial 42612)
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 30
$ws.Range("G4").Value = 2016

# Row 5 -> Transect 4, Date Sampled 9/1/2016 (serial 42614)
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2016

# Match the author's final selection/cursor position recorded in the diff.
$ws.Range("G5").Select() | Out-Null
